$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDF")

# Row 2: replace the old "abcd" / "Password is required" pair with a single
# "Sauce Labs Onesie" value (custom formatted) and drop column B's content.
$ws.Range("A2").Value = "Sauce Labs Onesie"
$ws.Range("B2").ClearContents() | Out-Null

# Apply the new font formatting (Calibri Light / custom dark color / major theme)
# to the updated cell.
$ws.Range("A2").Font.Name = "Calibri Light"
$ws.Range("A2").Font.Color = 1973527
$ws.Range("A2").Font.ThemeFont = 2

# New rows of numeric data below it.
$ws.Range("A3").Value = 7.99
$ws.Range("A4").Value = 6
$ws.Range("A5").Value = 129.94

# Column A widened to fit the new, longer text.
$ws.Columns.Item(1).ColumnWidth = 15.54296875

# This sheet (DDF) becomes the active / selected sheet, with A9 selected.
$ws.Activate() | Out-Null
$ws.Range("A9").Select() | Out-Null
